# Daily attendance processing
#
# For every row in the "Recorded By" column (G), if the cell holds a
# comma-separated list of recorders (e.g. "dnasr281@gmail.com, System"),
# rotate the list left by one position - i.e. move the first entry to the
# end of the list (e.g. "System, dnasr281@gmail.com"). Single-value cells
# (no comma) are left untouched, as are empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$col = 7  # Column G - "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $raw = $cell.Value2

    if ($raw -ne $null -and $raw -ne "") {
        $text = [string]$raw

        if ($text.Contains(",")) {
            $parts = $text.Split(",")
            $trimmed = @()
            foreach ($part in $parts) {
                $trimmed += $part.Trim()
            }

            if ($trimmed.Count -gt 1) {
                $rotated = $trimmed[1..($trimmed.Count - 1)] + $trimmed[0]
                $newValue = [string]::Join(", ", $rotated)
                $cell.Value = $newValue
            }
        }
    }
}
